# Update runner's complexity counter and visualizer's reduction calculation logics
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row -> new "Level of Complexity (LoC)" value (column E)
$updates = @{
    2  = 2
    4  = 2
    6  = 3
    7  = 5
    8  = 5
    9  = 7
    10 = 8
    11 = 9
    12 = 9
    13 = 10
    14 = 5
    15 = 5
    16 = 6
    17 = 6
    18 = 7
    19 = 7
    20 = 10
    21 = 10
    22 = 8
    23 = 8
    24 = 13
    25 = 13
    26 = 10
    27 = 10
    28 = 11
    29 = 11
    30 = 6
    31 = 6
    32 = 5
    33 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
